$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("Q5").Value = "$43.36"
$ws.Range("R5").Value = "FAIL"
